# Updated cryptos list on Tue Aug  6 21:38:37 UTC 2024 with GitHub Actions
# Refresh of Price (D) / Volume(1h) (E) columns; row 47/48 coin order swapped
# (Bittensor now ranks above WhiteBITCoin) with their refreshed data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D hold price strings that look numeric (e.g. "488.99", "5.82").
# Force those cells to Text format before assigning so Excel keeps the
# exact literal text (trailing zeros, precision) instead of silently
# coercing the string to a floating point number.

$ws.Range('D2').Value = '56.229.07'
$ws.Range('E2').Value = '  +3.21%  '
$ws.Range('D3').Value = '2.480.30'
$ws.Range('E3').Value = '  +1.43%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '488.99'
$ws.Range('E5').Value = '  +4.84%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '147.25'
$ws.Range('E6').Value = '  +10.60%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('E8').Value = '  +3.28%  '
$ws.Range('D9').Value = '2.489.99'
$ws.Range('E9').Value = '  +1.97%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '5.82'
$ws.Range('E10').Value = '  +9.37%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0971'
$ws.Range('E11').Value = '  +1.95%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.333'
$ws.Range('E12').Value = '  +5.30%  '
$ws.Range('E13').Value = '  +1.59%  '
$ws.Range('D14').Value = '2.915.59'
$ws.Range('E14').Value = '  +1.67%  '
$ws.Range('D15').Value = '56.236.01'
$ws.Range('E15').Value = '  +3.19%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '21.15'
$ws.Range('E16').Value = '  +6.73%  '
$ws.Range('E17').Value = '  +2.28%  '
$ws.Range('D18').Value = '2.486.52'
$ws.Range('E18').Value = '  +2.06%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.53'
$ws.Range('E19').Value = '  +7.96%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.09'
$ws.Range('E20').Value = '  +6.39%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '318.58'
$ws.Range('E21').Value = '  +2.63%  '
$ws.Range('E22').Value = '  +0.00%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.81'
$ws.Range('E23').Value = '  +7.99%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '58.45'
$ws.Range('E24').Value = '  +3.82%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.413'
$ws.Range('E25').Value = '  +7.04%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.999'
$ws.Range('E26').Value = '  -0.68%  '
$ws.Range('E27').Value = '  +4.52%  '
$ws.Range('D28').Value = '2.582.14'
$ws.Range('E28').Value = '  +2.16%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.62'
$ws.Range('E29').Value = '  +7.11%  '
$ws.Range('D30').Value = '0.0₃0792'
$ws.Range('E30').Value = '  +10.07%  '
$ws.Range('E31').Value = '  +0.15%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '149.24'
$ws.Range('E32').Value = '  +1.58%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '18.23'
$ws.Range('E33').Value = '  +2.62%  '
$ws.Range('E34').Value = '  +4.96%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.21'
$ws.Range('E35').Value = '  +4.22%  '
$ws.Range('E36').Value = '  +8.22%  '
$ws.Range('E37').Value = '  +5.40%  '
$ws.Range('E38').Value = '  +6.60%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '34.15'
$ws.Range('E39').Value = '  +3.86%  '
$ws.Range('E40').Value = '  +8.31%  '
$ws.Range('E41').Value = '  -0.02%  '
$ws.Range('E42').Value = '  +5.89%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.606'
$ws.Range('E43').Value = '  +1.70%  '
$ws.Range('E44').Value = '  +7.29%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '4.77'
$ws.Range('E45').Value = '  +14.40%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0926'
$ws.Range('E46').Value = '  +4.96%  '
$ws.Range('B47').Value = 'Bittensor'
$ws.Range('C47').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '258.45'
$ws.Range('E47').Value = '  +12.50%  '
$ws.Range('B48').Value = 'WhiteBITCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '10.20'
$ws.Range('E48').Value = '  +1.30%  '
$ws.Range('E49').Value = '  +4.98%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '17.61'
$ws.Range('E50').Value = '  +6.46%  '
$ws.Range('D51').Value = '1.876.09'
$ws.Range('E51').Value = '  -3.03%  '
